# "add account and type_list to Ready to deploy"
#
# 1) The "更新时间" (last-updated) column L for every data row moves from
#    2021-06-05 to 2021-06-18.
# 2) Accounts at rank 94 / rank 95 swap: row 95 now shows "胖了哥甄选"'s
#    info (and its live-stream/product counts), row 96 now shows
#    "中天潮购APP"'s info (and its counts). The rank numbers in column A
#    (94 / 95) stay where they are.
#
# Plain `Range.Value = "<digits>"` / `= "<date-like text>"` gets silently
# re-typed by Excel into a Number/Date (changing the cell's stored type
# and, because the old "General" style has to be swapped for a
# date/number-formatted one, its style index too). To keep these cells
# textual - matching the original workbook, where every one of these is a
# shared string - the values are staged as formulas on scratch cells,
# copied, and pasted back with "paste values" (xlPasteValues = -4163),
# which bakes in the literal text without re-triggering autodetection and
# leaves the destination's existing style alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基础数据")
$xlPasteValues = -4163

# --- 1) bump the "updated" date column (L2:L201) ---------------------------
$ws.Range("AA1:AA200").Formula = '=TEXT(DATE(2021,6,18),"yyyy-mm-dd")'
$ws.Range("AA1:AA200").Copy()
$ws.Range("L2:L201").PasteSpecial($xlPasteValues)
$ws.Range("AA1:AA200").ClearContents()

# --- 2) swap the account rows for rank 94 / rank 95 -------------------------
# Text fields - no autodetection risk, plain assignment keeps them as
# shared strings already.
$ws.Range("B95").Value = "胖了哥甄选"
$ws.Range("C95").Value = "admin2017666"
$ws.Range("D95").Value = "//p26.douyinpic.com/img/tos-cn-avt-0015/0d62c0a0c09358f7ff449b0e39b2462a~c5_1080x1080.webp?from=2956013662"
$ws.Range("F95").Value = "深圳市胖了哥甄选科技有限公司"
$ws.Range("G95").Value = "42.43w"

$ws.Range("B96").Value = "中天潮购APP"
$ws.Range("C96").Value = "wojiaolizai"
$ws.Range("D96").Value = "//p11.douyinpic.com/aweme/1080x1080/317a6000c592d64850617.heic?from=2956013662"
$ws.Range("F96").Value = "中天潮购科技官方账号"
$ws.Range("G96").Value = "151.51w"

# Numeric-looking fields (live-stream count / product count) - stage
# through formulas so they land as text, like the rest of the sheet.
$ws.Range("AA95").Formula = '=TEXT(33,"0")'
$ws.Range("AB95").Formula = '=TEXT(35,"0")'
$ws.Range("AA96").Formula = '=TEXT(40,"0")'
$ws.Range("AB96").Formula = '=TEXT(98,"0")'
$ws.Range("AA95:AB96").Copy()
$ws.Range("H95:I96").PasteSpecial($xlPasteValues)
$ws.Range("AA95:AB96").ClearContents()

$excel.CutCopyMode = 0
